# Update cryptocurrency price/volume figures (and reorder three rows)
# to match the refreshed coinranking.com snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.012.88"
$ws.Range("E2").Value = "  +1.65%  "
$ws.Range("D3").Value = "3.478.07"
$ws.Range("E3").Value = "  +1.99%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "416.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.58%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.630"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.05%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.732"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("E10").Value = "  +6.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.53"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.86"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.79%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000229"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.35%  "
$ws.Range("D14").Value = "4.034.38"
$ws.Range("E14").Value = "  +2.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.141"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.59"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.66%  "
$ws.Range("D17").Value = "3.488.33"
$ws.Range("E17").Value = "  +2.18%  "
$ws.Range("E18").Value = "  +0.77%  "
$ws.Range("E19").Value = "  -0.85%  "
$ws.Range("D20").Value = "62.848.06"
$ws.Range("E20").Value = "  +1.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "467.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "90.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.79%  "
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.72"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +13.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "33.61"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.33%  "
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.81%  "
$ws.Range("E30").Value = "  +1.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.20%  "
$ws.Range("E32").Value = "  -0.67%  "
$ws.Range("E33").Value = "  -1.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "41.11"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("E36").Value = "  +7.79%  "
$ws.Range("E37").Value = "  -3.06%  "
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("E39").Value = "  +2.94%  "
$ws.Range("E40").Value = "  +7.47%  "
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.324"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.42%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.135"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.72%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "148.16"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.89%  "
$ws.Range("E44").Value = "  -2.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.39%  "
$ws.Range("E46").Value = "  +2.65%  "
$ws.Range("D47").Value = "0.0₃0587"
$ws.Range("E47").Value = "  +37.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +9.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.64%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.145"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.52%  "
